# Generate Report for Handoff
#
# The localization-status report is regenerated: four files that were
# "Ready for handoff" (Priority = low) just got handed off, so their
# Priority flips to "ht" and their Latest Handoff Datetime is refreshed.
# The Overview sheet's "Latest HO Xliff Generate Date" for those same rows
# shares its text with the de-de sheet's handoff timestamp, so it moves
# together with it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    # zh-cn: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-09-04 18:35:45"

    # de-de: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-09-04 18:35:50"

    # Overview: Latest HO Xliff Generate Date (G) mirrors the de-de
    # handoff timestamp for these rows
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-04 18:35:50"
}
